$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.556.52"
$ws.Range("E2").Value = "  +2.62%  "
$ws.Range("D3").Value = "1.786.80"
$ws.Range("E3").Value = "  +1.08%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.39"
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.558"
$ws.Range("E6").Value = "  +1.39%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "33.11"
$ws.Range("E8").Value = "  +9.44%  "
$ws.Range("E9").Value = "  +2.04%  "
$ws.Range("E10").Value = "  +3.67%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0936"
$ws.Range("E11").Value = "  +1.43%  "
$ws.Range("D12").Value = "2.043.85"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.06"
$ws.Range("E13").Value = "  +11.45%  "
$ws.Range("D14").Value = "1.787.05"
$ws.Range("E14").Value = "  +1.15%  "
$ws.Range("E15").Value = "  +1.33%  "
$ws.Range("D16").Value = "34.543.83"
$ws.Range("E16").Value = "  +2.84%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.29"
$ws.Range("E17").Value = "  +2.84%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.55"
$ws.Range("E18").Value = "  +0.62%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "253.61"
$ws.Range("E19").Value = "  +1.38%  "
$ws.Range("D20").Value = "0.0₃0777"
$ws.Range("E20").Value = "  +5.88%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.00"
$ws.Range("E21").Value = "  -0.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.43"
$ws.Range("E22").Value = "  +2.28%  "
$ws.Range("E23").Value = "  +1.60%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.15"
$ws.Range("E24").Value = "  +0.68%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.41"
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.34"
$ws.Range("E26").Value = "  -0.21%  "
$ws.Range("E27").Value = "  +3.16%  "
$ws.Range("E28").Value = "  +0.41%  "
$ws.Range("E29").Value = "  -0.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.75"
$ws.Range("E30").Value = "  -0.65%  "
$ws.Range("E31").Value = "  +1.36%  "
$ws.Range("E32").Value = "  +0.51%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.59"
$ws.Range("E33").Value = "  +1.73%  "
$ws.Range("E34").Value = "  +4.46%  "
$ws.Range("D35").Value = "1.445.43"
$ws.Range("E35").Value = "  -1.93%  "
$ws.Range("E36").Value = "  -0.13%  "
$ws.Range("E37").Value = "  +3.35%  "
$ws.Range("E38").Value = "  -0.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "83.24"
$ws.Range("E39").Value = "  +0.46%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.81"
$ws.Range("E40").Value = "  +4.75%  "
$ws.Range("E41").Value = "  -0.27%  "
$ws.Range("E42").Value = "  +1.73%  "
$ws.Range("E43").Value = "  +0.32%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0504"
$ws.Range("E44").Value = "  -0.82%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.90"
$ws.Range("E45").Value = "  +2.87%  "
$ws.Range("E46").Value = "  -2.24%  "
$ws.Range("D47").Value = "1.941.05"
$ws.Range("E47").Value = "  +1.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "104.21"
$ws.Range("E48").Value = "  +7.38%  "
$ws.Range("B49").Value = "PaxDollar"
$ws.Range("C49").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.00"
$ws.Range("E49").Value = "  -0.20%  "
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.99"
$ws.Range("E50").Value = "  +1.37%  "
$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "49.26"
$ws.Range("E51").Value = "  -1.97%  "
